# session-06 slide-10: split the inline checklist into its own lines so
# each "checkbox" item renders on its own row instead of one long line.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$shape = $s.Shapes.Item(4)          # "Text 3" - the body placeholder with the numbered list
$tr = $shape.TextFrame.TextRange

$CR  = [char]13
$BOX = [char]0x2610                 # ☐

# --- Rebuild the whole paragraph list -------------------------------------
# Paragraph 7 ("3.  Checklist - ...") is replaced by 4 paragraphs: the
# "3.  Checklist:" lead-in plus one paragraph per checkbox item. Every other
# paragraph keeps its original text. Re-assigning the full TextRange.Text is
# what actually creates new paragraphs in this engine (per-paragraph /
# per-character sub-range assignment does not split text into paragraphs).
$lines = @(
    "A good template has three parts:",
    "",
    "1.  Title format $([char]0x2014) e.g., `"Equipment Issue: [Name] $([char]0x2014) [Problem]`"",
    "",
    "2.  Guided sections $([char]0x2014) Prompts: What happened? Where? How urgent?",
    "",
    "3.  Checklist:",
    "     $BOX Photos attached",
    "     $BOX Parts identified",
    "     $BOX Assigned",
    "",
    "Tip: Write your template as if explaining to someone new on the farm.",
    "What information would they need to give you?"
)

$tr.Text = [string]::Join($CR, $lines)

# --- Restore per-paragraph formatting --------------------------------------
# Re-assigning TextRange.Text above stamped every new paragraph with the
# formatting of paragraph 1 (sz=1800, bold, color 1E5128), so reapply the
# correct size / bold / italic / color for each paragraph.

function Set-ParaFormat($paraIndex, $size, $bold, $italic, $colorRgb) {
    # Re-stamping TextRange.Text only ever leaks bold=true / sz=1800 /
    # green from paragraph 1 - italic is already false for every new
    # paragraph, so only touch Font.Italic when a paragraph actually
    # needs to become italic (keeps the XML free of a redundant i="0").
    $para = $tr.Paragraphs($paraIndex, 1)
    $para.Font.Size = $size
    $para.Font.Bold = $bold
    if ($italic) {
        $para.Font.Italic = $italic
    }
    $para.Font.Color.RGB = $colorRgb
}

$DARKGRAY = 0x2C2C2C   # palindromic, no BGR swap needed
$GREEN    = 0x28511E   # BGR-swapped form of 1E5128

# 1: "A good template has three parts:" - already correct (sz 1800, bold, green)
# 2: blank - already correct

Set-ParaFormat 3 14 $false $false $DARKGRAY   # "1.  Title format ..."
# 4: blank - already correct

Set-ParaFormat 5 14 $false $false $DARKGRAY   # "2.  Guided sections ..."
# 6: blank - already correct

Set-ParaFormat 7  14 $false $false $DARKGRAY   # "3.  Checklist:"
Set-ParaFormat 8  13 $false $false $DARKGRAY   # "     ☐ Photos attached"
Set-ParaFormat 9  13 $false $false $DARKGRAY   # "     ☐ Parts identified"
Set-ParaFormat 10 13 $false $false $DARKGRAY   # "     ☐ Assigned"
# 11: blank - already correct

Set-ParaFormat 12 13 $false $true $GREEN   # "Tip: Write your template ..."
Set-ParaFormat 13 13 $false $true $GREEN   # "What information would they need to give you?"
